# "With Synchronization and Output"
# Adds a second output column (CartItems_output) next to the existing
# p_Quantity column on the "Global" sheet, mirroring the existing
# box-border formatting, and moves the active selection back to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")
$originalActiveSheet = $wb.ActiveSheet.Name

# --- New header cell (B1) -------------------------------------------------
$ws.Range("B1").Value = "CartItems_output"

# --- New data cell (B2) ----------------------------------------------------
$ws.Range("B2").Value = 2

# --- Give column B the same "boxed" look the old column A had --------------
# (right/top/bottom thin black border) on every data row.
foreach ($addr in "B2","B3","B4") {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(10).LineStyle = 1   # xlEdgeRight
  $rng.Borders.Item(10).Color = 0
  $rng.Borders.Item(8).LineStyle = 1    # xlEdgeTop
  $rng.Borders.Item(8).Color = 0
  $rng.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
  $rng.Borders.Item(9).Color = 0
}

# --- Column A no longer is the right edge of the box, drop its right border,
#     keeping only the top/bottom thin black border. ------------------------
foreach ($addr in "A2","A3","A4") {
  $rng = $ws.Range($addr)
  $rng.Borders.Item(10).LineStyle = -4142   # xlLineStyleNone
}

# --- Size column B to fit its new header text -------------------------------
$ws.Columns.Item(2).ColumnWidth = 14.93

# --- Move the selection back up to A2 on the Global sheet, then restore
#     whichever sheet/tab was active before this script ran. ---------------
$ws.Range("A2").Select()
$wb.Worksheets.Item($originalActiveSheet).Select()

Write-Output "done"
